$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D2:D51) to Text format so numeric-looking values
# (e.g. "242.59") are stored as literal strings, matching the source data,
# not auto-converted to numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.228.12"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.843.19"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "242.59"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "0.6637"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "44.84"
$ws.Range("E8").Value = "  +7.23%  "
$ws.Range("D9").Value = "0.07455"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "0.2956"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").Value = "23.45"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "0.07775"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "1.870.24"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "5.027"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "0.6733"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "83.64"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Value = "6.191"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "0.000008611"
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").Value = "29.245.72"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "2.113.46"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "227.68"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "12.57"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "7.186"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "1.000"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "158.77"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "0.1409"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "8.644"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "18.07"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "1.512"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "4.138"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "4.060"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "1.192"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "0.05342"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "1.880"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "0.7468"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").Value = "1.158"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").Value = "2.651"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "1.321.58"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "0.01801"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "2.756"
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").Value = "6.415"
$ws.Range("E42").Value = "  +7.36%  "
$ws.Range("D43").Value = "0.9047"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "103.31"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "2.001.61"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "66.01"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").Value = "0.00000000122"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").Value = "0.5141"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").Value = "0.07660"
$ws.Range("E50").Value = "  -8.42%  "
$ws.Range("D51").Value = "1.756"
$ws.Range("E51").Value = "  -0.51%  "

# Restore default (Normal/General) style on the Price column so the
# underlying cell style index is unchanged from before the edit.
$ws.Range("D2:D51").Style = "Normal"
